$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(55, 388, 7201, 1730, 452, 16, 2224),
    @(56, 387, 7588, 1815, 469, 17, 2314),
    @(57, 393, 7981, 1887, 482, 13, 2378),
    @(58, 495, 8476, 1945, 503, 21, 2416),
    @(59, 488, 8964, 2002, 514, 11, 2476),
    @(60, 436, 9400, 2075, 525, 11, 2556),
    @(61, 346, 9746, 2172, 533, 8, 2655),
    @(62, 347, 10093, 2326, 544, 11, 2811),
    @(63, 338, 10431, 2486, 556, 12, 2980)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

$excel.ActiveWindow.ScrollRow = 47
$ws.Range("C63").Select()
